$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 12969591
$ws.Range("I116").Value = 16668991
$ws.Range("J116").Value = 11912619
$ws.Range("K116").Value = 16668991
$ws.Range("L116").Value = 11912619
$ws.Range("M116").Value = -16665549
$ws.Range("N116").Value = -11919503

$ws.Range("H132").Value = 2526532.5
$ws.Range("I132").Value = 1362.825
$ws.Range("J132").Value = 27778228
$ws.Range("K132").Value = 4088.475
$ws.Range("L132").Value = 83334684
$ws.Range("M132").Value = -1558.475
$ws.Range("N132").Value = -83339744

$ws.Range("H137").Value = 26804106
$ws.Range("I137").Value = 5953332
$ws.Range("J137").Value = 89356424
$ws.Range("K137").Value = 17859996
$ws.Range("L137").Value = 268069272
$ws.Range("M137").Value = -17857446
$ws.Range("N137").Value = -268074372

$ws.Range("H141").Value = 1735
$ws.Range("I141").Value = 1735
$ws.Range("K141").Value = 5205
$ws.Range("M141").Value = -25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4188011.5
$ws.Range("I61").Value = 2193970.5
$ws.Range("K61").Value = 2193970.5
$ws.Range("M61").Value = -2193758.5

$ws.Range("H88").Value = 5833.222
$ws.Range("I88").Value = 2500
$ws.Range("J88").Value = 6785.5713
$ws.Range("K88").Value = 2500
$ws.Range("L88").Value = 6785.5713
$ws.Range("M88").Value = -2094
$ws.Range("N88").Value = -7597.5713

$ws.Range("H91").Value = 5833.222
$ws.Range("I91").Value = 2500
$ws.Range("J91").Value = 6785.5713
$ws.Range("K91").Value = 2500
$ws.Range("L91").Value = 6785.5713
$ws.Range("M91").Value = -1096
$ws.Range("N91").Value = -9593.5713

$ws.Range("H132").Value = 13414855
$ws.Range("I132").Value = 14498289
$ws.Range("K132").Value = 43494867
$ws.Range("M132").Value = -43492337

$ws.Range("H136").Value = 4188011.5
$ws.Range("I136").Value = 2193970.5
$ws.Range("K136").Value = 6581911.5
$ws.Range("M136").Value = -6579361.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1885.76
$ws.Range("I86").Value = 1914.7084
$ws.Range("J86").Value = 1191
$ws.Range("K86").Value = 1914.7084
$ws.Range("L86").Value = 1191
$ws.Range("M86").Value = -791.7084
$ws.Range("N86").Value = -3437

$ws.Range("H89").Value = 1885.76
$ws.Range("I89").Value = 1914.7084
$ws.Range("J89").Value = 1191
$ws.Range("K89").Value = 9573.541999999999
$ws.Range("L89").Value = 5955
$ws.Range("M89").Value = -3957.541999999999
$ws.Range("N89").Value = -17187

$ws.Range("H134").Value = 66967164
$ws.Range("I134").Value = 125002630
$ws.Range("K134").Value = 375007890
$ws.Range("M134").Value = -375005355

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1085604.4
$ws.Range("I134").Value = 4003.4243
$ws.Range("K134").Value = 12010.2729
$ws.Range("M134").Value = -9475.2729

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 436.2
$ws.Range("I13").Value = 81
$ws.Range("J13").Value = 525
$ws.Range("K13").Value = 243
$ws.Range("L13").Value = 1575
$ws.Range("M13").Value = -75
$ws.Range("N13").Value = -1911

$ws.Range("H70").Value = 6000
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

$ws.Range("H73").Value = 6000
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

$ws.Range("H76").Value = 5001.364
$ws.Range("J76").Value = 5001.5
$ws.Range("L76").Value = 15004.5
$ws.Range("N76").Value = -15770.5

$ws.Range("H79").Value = 5001.364
$ws.Range("J79").Value = 5001.5
$ws.Range("L79").Value = 15004.5
$ws.Range("N79").Value = -17656.5

$ws.Range("H82").Value = 22808
$ws.Range("I82").Value = 480
$ws.Range("J82").Value = 25288.889
$ws.Range("K82").Value = 1440
$ws.Range("L82").Value = 75866.667
$ws.Range("M82").Value = -1034
$ws.Range("N82").Value = -76678.667

$ws.Range("H85").Value = 22808
$ws.Range("I85").Value = 480
$ws.Range("J85").Value = 25288.889
$ws.Range("K85").Value = 1440
$ws.Range("L85").Value = 75866.667
$ws.Range("M85").Value = -36
$ws.Range("N85").Value = -78674.667

$ws.Range("H94").Value = 2505.25
$ws.Range("I94").Value = 1664.6666
$ws.Range("J94").Value = 5027
$ws.Range("K94").Value = 4993.9998
$ws.Range("L94").Value = 15081
$ws.Range("M94").Value = -4317.9998
$ws.Range("N94").Value = -16433

$ws.Range("H97").Value = 1854.2222
$ws.Range("I97").Value = 897.5
$ws.Range("J97").Value = 2619.6
$ws.Range("K97").Value = 2692.5
$ws.Range("L97").Value = 7858.799999999999
$ws.Range("M97").Value = -2196.5
$ws.Range("N97").Value = -8850.799999999999

$ws.Range("H106").Value = 5989.875
$ws.Range("J106").Value = 5989.875
$ws.Range("L106").Value = 17969.625
$ws.Range("N106").Value = -19861.625

$ws.Range("H109").Value = 3629.1892
$ws.Range("I109").Value = 900.1429000000001
$ws.Range("J109").Value = 4265.967
$ws.Range("K109").Value = 2700.4287
$ws.Range("L109").Value = 12797.901
$ws.Range("M109").Value = -1660.4287
$ws.Range("N109").Value = -14877.901

$ws.Range("H112").Value = 2390
$ws.Range("I112").Value = 1015
$ws.Range("J112").Value = 6515
$ws.Range("K112").Value = 3045
$ws.Range("L112").Value = 19545
$ws.Range("M112").Value = -1937
$ws.Range("N112").Value = -21761

$ws.Range("H115").Value = 3661.0908
$ws.Range("I115").Value = 950
$ws.Range("J115").Value = 3932.2
$ws.Range("K115").Value = 2850
$ws.Range("L115").Value = 11796.6
$ws.Range("M115").Value = -1675
$ws.Range("N115").Value = -14146.6

$ws.Range("H118").Value = 3670.6667
$ws.Range("I118").Value = 2990
$ws.Range("K118").Value = 8970
$ws.Range("M118").Value = -7727

$ws.Range("H121").Value = 2153607.2
$ws.Range("I121").Value = 325
$ws.Range("J121").Value = 2302109.5
$ws.Range("K121").Value = 975
$ws.Range("L121").Value = 6906328.5
$ws.Range("M121").Value = 335
$ws.Range("N121").Value = -6908948.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 30777248
$ws.Range("J132").Value = 15154448
$ws.Range("L132").Value = 45463344
$ws.Range("N132").Value = -45468404

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 67095.2
$ws.Range("I16").Value = 77369.84
$ws.Range("K16").Value = 77369.84
$ws.Range("M16").Value = -77199.84

$ws.Range("H132").Value = 4178091.5
$ws.Range("I132").Value = 6678346.5
$ws.Range("J132").Value = 10999.667
$ws.Range("K132").Value = 20035039.5
$ws.Range("L132").Value = 32999.001
$ws.Range("M132").Value = -20032509.5
$ws.Range("N132").Value = -38059.001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 21555
$ws.Range("I107").Value = 18251.818
$ws.Range("J107").Value = 33666.668
$ws.Range("K107").Value = 54755.454
$ws.Range("L107").Value = 101000.004
$ws.Range("M107").Value = -52835.454
$ws.Range("N107").Value = -104840.004

$ws.Range("H132").Value = 1599864
$ws.Range("I132").Value = 5625.6665
$ws.Range("J132").Value = 4332844
$ws.Range("K132").Value = 16876.9995
$ws.Range("L132").Value = 12998532
$ws.Range("M132").Value = -14346.9995
$ws.Range("N132").Value = -13003592
